# Auto-generated script applying the cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (column D) and Volume(1h) (column E) cells with the refreshed snapshot values.
# Column D values that look like plain numbers are entered with a leading apostrophe and the
# cell style is reset to Normal afterwards so they stay plain text (matching the source data)
# without leaving a numeric format override on the cell.
$ws.Range("D2").Value = "26.851.16"
$ws.Range("D3").Value = "1.802.00"
$ws.Range("E3").Value = "  -1.32%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'309.26"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Range("E5").Value = "  -1.82%  "
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").Value = "'0.4661"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Range("E7").Value = "  +4.06%  "
$ws.Range("D8").Value = "'0.3695"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Range("E8").Value = "  -2.04%  "
$ws.Range("D9").Value = "'0.07362"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Range("E9").Value = "  -1.30%  "
$ws.Range("D10").Value = "'0.8689"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Range("E10").Value = "  -2.13%  "
$ws.Range("D11").Value = "'20.35"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Range("E11").Value = "  -3.20%  "
$ws.Range("D12").Value = "1.822.86"
$ws.Range("E12").Value = "  +0.18%  "
$ws.Range("D13").Value = "'5.355"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Range("E13").Value = "  -1.99%  "
$ws.Range("D14").Value = "'92.71"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Range("E14").Value = "  -1.44%  "
$ws.Range("D15").Value = "'6.503"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Range("E15").Value = "  -3.74%  "
$ws.Range("D16").Value = "'0.07026"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Range("E16").Value = "  -1.34%  "
$ws.Range("D17").Value = "'1.000"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Range("E17").Value = "  -0.11%  "
$ws.Range("D18").Value = "'0.000008685"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Range("E18").Value = "  -1.23%  "
$ws.Range("E19").Value = "  -0.01%  "
$ws.Range("E20").Value = "  -3.20%  "
$ws.Range("D21").Value = "26.846.18"
$ws.Range("D22").Value = "'5.280"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Range("E22").Value = "  -2.45%  "
$ws.Range("E23").Value = "  -3.69%  "
$ws.Range("D24").Value = "2.011.70"
$ws.Range("E24").Value = "  -2.09%  "
$ws.Range("D25").Value = "'1.904"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Range("E25").Value = "  -2.97%  "
$ws.Range("D26").Value = "'151.71"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Range("E26").Value = "  +0.15%  "
$ws.Range("D27").Value = "'18.33"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Range("E27").Value = "  -1.83%  "
$ws.Range("D28").Value = "'2.133"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Range("E28").Value = "  -8.33%  "
$ws.Range("D29").Value = "'5.244"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Range("E29").Value = "  -3.14%  "
$ws.Range("D30").Value = "'116.10"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Range("E30").Value = "  -1.61%  "
$ws.Range("D31").Value = "'0.08920"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Range("E31").Value = "  +0.31%  "
$ws.Range("D32").Value = "'0.7598"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Range("E32").Value = "  -4.10%  "
$ws.Range("D33").Value = "'2.940"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Range("E33").Value = "  +0.58%  "
$ws.Range("E34").Value = "  -4.58%  "
$ws.Range("D35").Value = "'4.460"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Range("E35").Value = "  -3.10%  "
$ws.Range("D36").Value = "'1.0000"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("E37").Value = "  -0.91%  "
$ws.Range("D38").Value = "'0.01952"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Range("E38").Value = "  -1.85%  "
$ws.Range("D39").Value = "'0.05246"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Range("E39").Value = "  -1.26%  "
$ws.Range("E40").Value = "  +1.71%  "
$ws.Range("D41").Value = "'7.208"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Range("E41").Value = "  -1.46%  "
$ws.Range("E44").Value = "  -3.57%  "
$ws.Range("D45").Value = "'8.499"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Range("E45").Value = "  -2.11%  "
$ws.Range("D46").Value = "'0.5010"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Range("E46").Value = "  -1.99%  "
$ws.Range("D47").Value = "'10.26"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Range("E47").Value = "  -4.16%  "
$ws.Range("D48").Value = "'103.92"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Range("E48").Value = "  -1.45%  "
$ws.Range("D49").Value = "'0.9998"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Range("E49").Value = "  -0.04%  "
$ws.Range("E50").Value = "  -2.34%  "
$ws.Range("D51").Value = "'0.06286"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Range("E51").Value = "  -2.01%  "

# Rows 42 and 43 swap rank order: RenderToken now ranks above TheSandbox, with refreshed price/volume
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").Value = "'2.366"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Range("E42").Value = "  +2.83%  "

$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").Value = "'0.5294"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Range("E43").Value = "  -1.13%  "
